# Translate the English item texts (column B, rows 2-17) of the
# "Tobler et al., 2012" scale workbook into German.
# Column A (the item codes like "01ccc", "02fop", ...) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$translations = @{
    2  = "Wir müssen das empfindliche Gleichgewicht des Klimas schützen."
    3  = "Klimaschutz ist wichtig für unsere Zukunft."
    4  = "Ich mache mir Sorgen über den Zustand des Klimas."
    5  = "Der Klimawandel hat schwerwiegende Folgen für Mensch und Natur."
    6  = "Klimaschutzmaßnahmen werden von wenigen mächtigen Personen bestimmt, ich als einzelner Bürger habe keinen Einfluss."
    7  = "Mit meinem Verhalten kann ich das Klima nicht beeinflussen, denn es liegt in den Händen der Industrie."
    8  = "Als normaler Bürger kann ich die Entscheidungen der Regierung zum Klimaschutz beeinflussen."
    9  = "Ich fühle mich in der Lage, einen Beitrag zum Klimaschutz zu leisten."
    10 = "Wenn ich versuchen würde, mich klimafreundlich zu verhalten, würde sich das sicher positiv auf das Klima auswirken."
    11 = "Der Klimawandel und seine Folgen werden in den Medien übertrieben dargestellt."
    12 = "Der Klimawandel ist eine Masche."
    13 = "Solange Meteorologen nicht einmal in der Lage sind, das Wetter genau vorherzusagen, lässt sich auch das Klima nicht zuverlässig vorhersagen."
    14 = "Es gibt größere Probleme als Klimaschutz."
    15 = "Ich fühle mich durch den Klimawandel nicht bedroht."
    16 = "Die Auswirkungen des Klimawandels sind unvorhersehbar; daher ist mein klimafreundliches Verhalten sinnlos."
    17 = "Der Klimaschutz behindert unnötig das Wirtschaftswachstum."
}

foreach ($row in $translations.Keys) {
    $ws.Cells.Item($row, 2).Value = $translations[$row]
}

# Let Excel recompute the (word-wrapped) row heights for the edited rows,
# then pin down the exact heights Excel settled on for the new (longer)
# German text so the wrapped row heights match the saved workbook.
$ws.Rows.AutoFit() | Out-Null

$rowHeights = @{
    2  = 30
    6  = 45
    9  = 30
    10 = 45
    16 = 45
}
foreach ($row in $rowHeights.Keys) {
    $ws.Rows.Item($row).RowHeight = $rowHeights[$row]
}

# Reset the zoom back to the default 100% and move the active selection,
# matching the saved state captured in the workbook.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("F16").Select()
